{"js": "// Replace the math-problem text in each table cell according to the\n// mapping below. The mapping is applied in document order so that the\n// transient collision (765\u00f78= -> 189\u00f74=, while an earlier cell 189\u00f74=\n// -> 545\u00f79=) resolves correctly: the original \"189\u00f74=\" is consumed\n// before the new \"189\u00f74=\" is produced.\nconst replacements = [\n  [\"115\u00f77=\", \"538\u00f76=\"],\n  [\"563\u00f79=\", \"989\u00f76=\"],\n  [\"269\u00f75=\", \"710\u00f75=\"],\n  [\"794\u00f76=\", \"974\u00f75=\"],\n  [\"184\u00f76=\", \"757\u00f75=\"],\n  [\"113\u00f73=\", \"601\u00f75=\"],\n  [\"902\u00f72=\", \"524\u00f78=\"],\n  [\"553\u00f78=\", \"730\u00f75=\"],\n  [\"690\u00f73=\", \"942\u00f74=\"],\n  [\"474\u00f77=\", \"650\u00f74=\"],\n  [\"228\u00f72=\", \"393\u00f72=\"],\n  [\"189\u00f74=\", \"545\u00f79=\"],\n  [\"713\u00f75=\", \"416\u00f76=\"],\n  [\"490\u00f77=\", \"221\u00f76=\"],\n  [\"526\u00f73=\", \"718\u00f79=\"],\n  [\"281\u00f79=\", \"869\u00f73=\"],\n  [\"296\u00f78=\", \"825\u00f75=\"],\n  [\"505\u00f74=\", \"903\u00f72=\"],\n  [\"876\u00f72=\", \"250\u00f78=\"],\n  [\"588\u00f73=\", \"851\u00f75=\"],\n  [\"765\u00f78=\", \"189\u00f74=\"],\n  [\"653\u00f76=\", \"655\u00f73=\"],\n  [\"628\u00f76=\", \"420\u00f73=\"],\n  [\"276\u00f79=\", \"850\u00f76=\"],\n  [\"285\u00f78=\", \"629\u00f72=\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items,text,font\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${from}`);\n  }\n\n  // Replace only the first occurrence \u2014 each \"from\" string is unique in\n  // the document at the time it is processed.\n  results.items[0].insertText(to, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the math-problem text in each table cell according to the\n# mapping below. The mapping is applied in document order so that the\n# transient collision (765\u00f78= -> 189\u00f74=, while an earlier cell 189\u00f74=\n# -> 545\u00f79=) resolves correctly: the original \"189\u00f74=\" is consumed\n# before the new \"189\u00f74=\" is produced.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"115\u00f77=\", \"538\u00f76=\"),\n    @(\"563\u00f79=\", \"989\u00f76=\"),\n    @(\"269\u00f75=\", \"710\u00f75=\"),\n    @(\"794\u00f76=\", \"974\u00f75=\"),\n    @(\"184\u00f76=\", \"757\u00f75=\"),\n    @(\"113\u00f73=\", \"601\u00f75=\"),\n    @(\"902\u00f72=\", \"524\u00f78=\"),\n    @(\"553\u00f78=\", \"730\u00f75=\"),\n    @(\"690\u00f73=\", \"942\u00f74=\"),\n    @(\"474\u00f77=\", \"650\u00f74=\"),\n    @(\"228\u00f72=\", \"393\u00f72=\"),\n    @(\"189\u00f74=\", \"545\u00f79=\"),\n    @(\"713\u00f75=\", \"416\u00f76=\"),\n    @(\"490\u00f77=\", \"221\u00f76=\"),\n    @(\"526\u00f73=\", \"718\u00f79=\"),\n    @(\"281\u00f79=\", \"869\u00f73=\"),\n    @(\"296\u00f78=\", \"825\u00f75=\"),\n    @(\"505\u00f74=\", \"903\u00f72=\"),\n    @(\"876\u00f72=\", \"250\u00f78=\"),\n    @(\"588\u00f73=\", \"851\u00f75=\"),\n    @(\"765\u00f78=\", \"189\u00f74=\"),\n    @(\"653\u00f76=\", \"655\u00f73=\"),\n    @(\"628\u00f76=\", \"420\u00f73=\"),\n    @(\"276\u00f79=\", \"850\u00f76=\"),\n    @(\"285\u00f78=\", \"629\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.Text = $to\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 0\n    [void]$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, 0, $false, $to, 2)\n}\n"}
